$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates (column C) ---
$ws.Range("C2").Value = "Repository is null"
$ws.Range("C4").Value = "Repository"
$ws.Range("C5").Value = "Tag"

# --- New (currently-empty) column D formatting: red Times New Roman 12, wrap text ---
$ws.Columns.Item(4).Font.Name = "Times New Roman"
$ws.Columns.Item(4).Font.Size = 12
$ws.Columns.Item(4).Font.Color = 255
$ws.Columns.Item(4).WrapText = $true

# --- Normalize column C font color (drop inherited theme color) ---
$ws.Range("C1:C6").Font.Color = 0

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moves to C9 ---
$ws.Range("C9").Select()
